$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 36856.23
$ws.Range("I21").Value = 28395.223
$ws.Range("J21").Value = 55893.5
$ws.Range("K21").Value = 28395.223
$ws.Range("L21").Value = 55893.5
$ws.Range("M21").Value = -27927.223
$ws.Range("N21").Value = -56829.5
$ws.Range("H23").Value = 36856.23
$ws.Range("I23").Value = 28395.223
$ws.Range("J23").Value = 55893.5
$ws.Range("K23").Value = 28395.223
$ws.Range("L23").Value = 55893.5
$ws.Range("M23").Value = -28161.223
$ws.Range("N23").Value = -56361.5
$ws.Range("H38").Value = 469.85715
$ws.Range("I38").Value = 337.8
$ws.Range("J38").Value = 800
$ws.Range("K38").Value = 1013.4
$ws.Range("L38").Value = 2400
$ws.Range("M38").Value = -641.4000000000001
$ws.Range("N38").Value = -3144
$ws.Range("H58").Value = 1273.3334
$ws.Range("I58").Value = 220
$ws.Range("J58").Value = 1800
$ws.Range("K58").Value = 660
$ws.Range("L58").Value = 5400
$ws.Range("M58").Value = -510
$ws.Range("N58").Value = -5700
$ws.Range("H117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("N117").ClearContents()
$ws.Range("H138").Value = 2191.8396
$ws.Range("I138").Value = 1346.8857
$ws.Range("J138").Value = 2834.739
$ws.Range("K138").Value = 4040.6571
$ws.Range("L138").Value = 8504.217000000001
$ws.Range("M138").Value = 1099.3429
$ws.Range("N138").Value = -18784.217

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 12437.707
$ws.Range("I32").Value = 11176.849
$ws.Range("J32").Value = 25802.8
$ws.Range("K32").Value = 11176.849
$ws.Range("L32").Value = 25802.8
$ws.Range("M32").Value = -10889.849
$ws.Range("N32").Value = -26376.8
$ws.Range("H54").Value = 8191.875
$ws.Range("I54").Value = 5000
$ws.Range("J54").Value = 8647.857
$ws.Range("K54").Value = 5000
$ws.Range("L54").Value = 8647.857
$ws.Range("M54").Value = -4231
$ws.Range("N54").Value = -10185.857
$ws.Range("H101").Value = 38602
$ws.Range("J101").Value = 38602
$ws.Range("L101").Value = 38602
$ws.Range("N101").Value = -45092

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").ClearContents()
$ws.Range("H134").Value = 1863.1428
$ws.Range("I134").Value = 1463.6487
$ws.Range("K134").Value = 4390.9461
$ws.Range("M134").Value = -1855.9461

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H14").Value = 23502.75
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 23502.75
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 23502.75
$ws.Range("M14").ClearContents()
$ws.Range("N14").Value = -23842.75
$ws.Range("H31").Value = 3404661.8
$ws.Range("I31").Value = 1434.4694
$ws.Range("J31").Value = 6807889
$ws.Range("K31").Value = 1434.4694
$ws.Range("L31").Value = 6807889
$ws.Range("M31").Value = -1139.4694
$ws.Range("N31").Value = -6808479
$ws.Range("H34").Value = 3404661.8
$ws.Range("I34").Value = 1434.4694
$ws.Range("J34").Value = 6807889
$ws.Range("K34").Value = 1434.4694
$ws.Range("L34").Value = 6807889
$ws.Range("M34").Value = -1232.4694
$ws.Range("N34").Value = -6808293
$ws.Range("H122").Value = 48903.4
$ws.Range("I122").Value = 50899.375
$ws.Range("J122").Value = 1000
$ws.Range("K122").Value = 152698.125
$ws.Range("L122").Value = 3000
$ws.Range("M122").Value = -150248.125
$ws.Range("N122").Value = -7900
$ws.Range("H132").Value = 343197.78
$ws.Range("I132").Value = 1222.3125
$ws.Range("J132").Value = 1559110.5
$ws.Range("K132").Value = 3666.9375
$ws.Range("L132").Value = 4677331.5
$ws.Range("M132").Value = -1136.9375
$ws.Range("N132").Value = -4682391.5
$ws.Range("H134").Value = 760238.6
$ws.Range("I134").Value = 455348.97
$ws.Range("K134").Value = 1366046.91
$ws.Range("M134").Value = -1363511.91

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 310
$ws.Range("J17").Value = 450
$ws.Range("L17").Value = 1350
$ws.Range("N17").Value = -1688
$ws.Range("H34").Value = 1446.375
$ws.Range("J34").Value = 1676.3077
$ws.Range("L34").Value = 5028.9231
$ws.Range("N34").Value = -5196.9231
$ws.Range("H35").Value = 1000
$ws.Range("I35").Value = 1000
$ws.Range("J35").Value = 1000
$ws.Range("K35").Value = 3000
$ws.Range("L35").Value = 3000
$ws.Range("M35").Value = -2712
$ws.Range("N35").Value = -3576
$ws.Range("H39").Value = 710.44446
$ws.Range("J39").Value = 710.44446
$ws.Range("L39").Value = 2131.33338
$ws.Range("N39").Value = -2719.33338
$ws.Range("H55").Value = 1950
$ws.Range("J55").Value = 1950
$ws.Range("L55").Value = 5850
$ws.Range("N55").Value = -6204
$ws.Range("H68").Value = 1172.1704
$ws.Range("I68").Value = 786.1429000000001
$ws.Range("J68").Value = 1352.3167
$ws.Range("K68").Value = 2358.4287
$ws.Range("L68").Value = 4056.9501
$ws.Range("M68").Value = -1547.4287
$ws.Range("N68").Value = -5678.9501
$ws.Range("H71").Value = 1172.1704
$ws.Range("I71").Value = 786.1429000000001
$ws.Range("J71").Value = 1352.3167
$ws.Range("K71").Value = 7075.2861
$ws.Range("L71").Value = 12170.8503
$ws.Range("M71").Value = -3019.2861
$ws.Range("N71").Value = -20282.8503
$ws.Range("H98").Value = 368.33334
$ws.Range("I98").Value = 162.5
$ws.Range("K98").Value = 487.5
$ws.Range("M98").Value = 1010.5
$ws.Range("H105").Value = 336999.34
$ws.Range("J105").Value = 336999.34
$ws.Range("L105").Value = 1010998.02
$ws.Range("N105").Value = -1016240.02
$ws.Range("H113").Value = 2018.9117
$ws.Range("I113").Value = 2389.077
$ws.Range("J113").Value = 815.875
$ws.Range("K113").Value = 7167.231000000001
$ws.Range("L113").Value = 2447.625
$ws.Range("M113").Value = -4997.231000000001
$ws.Range("N113").Value = -6787.625
$ws.Range("H123").Value = 2181.111
$ws.Range("I123").Value = 2181.111
$ws.Range("J123").Value = 0
$ws.Range("K123").Value = 6543.333
$ws.Range("L123").Value = 0
$ws.Range("M123").Value = -4093.333
$ws.Range("N123").ClearContents()
$ws.Range("H126").Value = 4894.125
$ws.Range("J126").Value = 4894.125
$ws.Range("L126").Value = 14682.375
$ws.Range("N126").Value = -24562.375

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H42").Value = 30000
$ws.Range("J42").Value = 30000
$ws.Range("L42").Value = 30000
$ws.Range("N42").Value = -30970
$ws.Range("H107").Value = 3881.2
$ws.Range("I107").Value = 421.42856
$ws.Range("J107").Value = 11954
$ws.Range("K107").Value = 421.42856
$ws.Range("L107").Value = 11954
$ws.Range("M107").Value = 1498.57144
$ws.Range("N107").Value = -15794
$ws.Range("H115").Value = 30000
$ws.Range("J115").Value = 30000
$ws.Range("L115").Value = 30000
$ws.Range("N115").Value = -32350
$ws.Range("H134").Value = 14108.667
$ws.Range("J134").Value = 14108.667
$ws.Range("L134").Value = 42326.001
$ws.Range("N134").Value = -47396.001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 3642.1072
$ws.Range("I68").Value = 3624.1667
$ws.Range("J68").Value = 3655.5625
$ws.Range("K68").Value = 3624.1667
$ws.Range("L68").Value = 3655.5625
$ws.Range("M68").Value = -2875.1667
$ws.Range("N68").Value = -5153.5625
$ws.Range("H71").Value = 3642.1072
$ws.Range("I71").Value = 3624.1667
$ws.Range("J71").Value = 3655.5625
$ws.Range("K71").Value = 18120.8335
$ws.Range("L71").Value = 18277.8125
$ws.Range("M71").Value = -14376.8335
$ws.Range("N71").Value = -25765.8125

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 967660.25
$ws.Range("I132").Value = 1359671.2
$ws.Range("J132").Value = 2710.077
$ws.Range("K132").Value = 4079013.6
$ws.Range("L132").Value = 8130.231000000001
$ws.Range("M132").Value = -4076483.6
$ws.Range("N132").Value = -13190.231
